$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve existing formatting of A1 (header) and A2 before we touch column A ---
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122)   # xlPasteFormats

# Drop column A's stale auto-fit/leftover style bookkeeping (bestFit + default column style)
$ws.Columns.Item(1).ClearFormats()

# Restore A1 / A2 formatting (bold header border, bordered A2) that ClearFormats wiped
$ws.Range("F1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Clean up the temporary helper cells
$ws.Range("F1:F2").Clear()

# Widen column A to fit the new phone-number content (manual width, no longer "best fit")
$ws.Columns.Item(1).ColumnWidth = 20

# Load the phone number into A2 (keeps its existing bordered style)
$ws.Range("A2").Value = 932000076

# The rest of column A (A3:A10) loses its placeholder border now that a real row has data
$ws.Range("A3:A10").Borders.LineStyle = -4142   # xlLineStyleNone

# Move the active selection, matching where the user clicked next (phone cell shifted right)
$ws.Range("E6").Select()
